# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" column (D) for the row corresponding
# to the c738fd0b-...md file (row 5) on both the zh-cn and de-de sheets,
# recording the timestamps of the newly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-27 07:53:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-27 07:53:54"
